$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill existing rows 6 and 7 (columns D:G) with values matching rows 2/3 pattern ---
$ws.Range("D6").Value = 47.39
$ws.Range("E6").Value = 71.41
$ws.Range("F6").Value = 40.43
$ws.Range("G6").Value = 72.59

$ws.Range("D7").Value = 47.39
$ws.Range("E7").Value = 71.41
$ws.Range("F7").Value = 40.43
$ws.Range("G7").Value = 72.59

$ws.Range("D6:G7").Interior.Color = 65535

# --- Row 8: scheduler.step param set, div_factor, batch 100 ---
$ws.Range("A8").Value = 4
$ws.Range("B8").Value = "div_factor"
$ws.Range("C8").Value = 100
$ws.Range("D8").Value = 68.48
$ws.Range("E8").Value = 86.98
$ws.Range("F8").Value = 46.47
$ws.Range("G8").Value = 86.85

# --- Row 9: final_div_factor, batch 100 ---
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "final_div_factor"
$ws.Range("C9").Value = 100
$ws.Range("D9").Value = 68.48
$ws.Range("E9").Value = 86.98
$ws.Range("F9").Value = 46.47
$ws.Range("G9").Value = 86.85

# --- Row 10: div_factor, batch 500 (highlighted style, same as rows 2/3/6/7) ---
$ws.Range("A10").Value = 5
$ws.Range("B10").Value = "div_factor"
$ws.Range("C10").Value = 500
$ws.Range("D10").Value = 68.09
$ws.Range("E10").Value = 86.76
$ws.Range("F10").Value = 57.65
$ws.Range("G10").Value = 86.2

# --- Row 11: final_div_factor, batch 500 (highlighted style) ---
$ws.Range("A11").Value = 5
$ws.Range("B11").Value = "final_div_factor"
$ws.Range("C11").Value = 500
$ws.Range("D11").Value = 68.09
$ws.Range("E11").Value = 86.76
$ws.Range("F11").Value = 57.65
$ws.Range("G11").Value = 86.2

$ws.Range("A10:G11").Interior.Color = 65535

# --- Row 12: div_factor, batch 100 ---
$ws.Range("A12").Value = 6
$ws.Range("B12").Value = "div_factor"
$ws.Range("C12").Value = 100
$ws.Range("D12").Value = 68.26
$ws.Range("E12").Value = 87.12
$ws.Range("F12").Value = 30.35
$ws.Range("G12").Value = 86.74

# --- Row 13: final_div_factor, batch 500 ---
$ws.Range("A13").Value = 6
$ws.Range("B13").Value = "final_div_factor"
$ws.Range("C13").Value = 500
$ws.Range("D13").Value = 68.26
$ws.Range("E13").Value = 87.12
$ws.Range("F13").Value = 30.35
$ws.Range("G13").Value = 86.74

# --- Row 14: div_factor, batch 10000, highlighted red (new scheduler.step entry) ---
$ws.Range("A14").Value = 5
$ws.Range("B14").Value = "div_factor"
$ws.Range("C14").Value = 10000

# --- Row 15: final_div_factor, batch 10000, highlighted red ---
$ws.Range("A15").Value = 5
$ws.Range("B15").Value = "final_div_factor"
$ws.Range("C15").Value = 10000

$ws.Range("A14:C15").Interior.Color = 255  # RGB(255,0,0) red, encoded as R + G*256 + B*65536

$ws.Range("F22").Select()
